# Applies the "transition data" edit described in the commit:
#   "lots of transition data that may need filtering"
#
# Summary of changes on Sheet1:
#  - The old "GR" (Green roof) In/Out sample values in D3:E8 are replaced
#    with a fresh set of values (sourced "from StormTac").
#  - The old D3:E8 values are preserved in a new side table in N1:O9,
#    labeled "From US BMP database" (still the "GR" category, In/Out).
#  - A new comparison table is added in K1:L8 ("from StormTac": traditional
#    roof vs. green roof), mirroring the new D3:E8 values.
#  - A threaded comment is added on N4 (same note as the existing D3 one).
#  - Sheet1 becomes the active/selected sheet (selection on P3).
#
# NOTE: this runtime's `Range.Value` getter is unreliable (returns a
# descriptor string instead of the cell's value), so this script always
# reads/writes through `Value2` or uses literal constants instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) New small comparison table header: K1:L8 ("from StormTac")
#    Writing the brand-new strings in this particular order matches the
#    order they'd be appended to the shared-strings table.
# ---------------------------------------------------------------------
$ws.Range("K2").Value2 = "traditional roof"
$ws.Range("L2").Value2 = "Green roof"
$ws.Range("K1").Value2 = "from StormTac"
$ws.Range("N1").Value2 = "From US BMP database"

# ---------------------------------------------------------------------
# 2) Preserve the OLD "GR" In/Out values into the new N:O side table
#    before they get overwritten in D:E below.
# ---------------------------------------------------------------------
$ws.Range("N2").Value2 = "GR"
$ws.Range("N3").Value2 = "In "
$ws.Range("O3").Value2 = "Out"

$ws.Range("N4").Value2 = 0
$ws.Range("O4").Value2 = 4
$ws.Range("O5").Value2 = 0.42
$ws.Range("O6").Value2 = 1.866
$ws.Range("O7").Value2 = 0.0228
$ws.Range("O8").Value2 = 0.00067
$ws.Range("O9").Value2 = 0.018

# ---------------------------------------------------------------------
# 3) New "Green roof" values replacing the old D:E column
# ---------------------------------------------------------------------
$ws.Range("D3").Value2 = 22
$ws.Range("E3").Value2 = 19

$ws.Range("D4").Value2 = 0.053
$ws.Range("E4").Value2 = 0.059

$ws.Range("D5").Value2 = 1.7
$ws.Range("E5").Value2 = 1.8

$ws.Range("D6").Value2 = 0.022
$ws.Range("E6").Value2 = 0.016

$ws.Range("D7").Value2 = 0.005
$ws.Range("E7").Value2 = 0.001

$ws.Range("D8").Value2 = 0.08
$ws.Range("E8").Value2 = 0.023

# Mirror the new Green-roof values into the K:L comparison table.
$ws.Range("K3").Value2 = 22
$ws.Range("L3").Value2 = 19

$ws.Range("K4").Value2 = 0.053
$ws.Range("L4").Value2 = 0.059

$ws.Range("K5").Value2 = 1.7
$ws.Range("L5").Value2 = 1.8

$ws.Range("K6").Value2 = 0.022
$ws.Range("L6").Value2 = 0.016

$ws.Range("K7").Value2 = 0.005
$ws.Range("L7").Value2 = 0.001

$ws.Range("K8").Value2 = 0.08
$ws.Range("L8").Value2 = 0.023

# ---------------------------------------------------------------------
# 4) Formatting
# ---------------------------------------------------------------------

# K2:L2 header cells -> bold 12pt look (same as the other table headers)
# plus a thin left border to separate the two new columns.
$headerRange = $ws.Range("K2:L2")
$headerRange.Font.Name = "Aptos Narrow"
$headerRange.Font.Size = 12
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# N2:O2 header ("GR") uses the same bold look as the other merged headers.
$ws.Range("N2:O2").Font.Name = "Aptos Narrow"
$ws.Range("N2:O2").Font.Size = 12
$ws.Range("N2:O2").Font.Bold = $true
$ws.Range("N2:O2").HorizontalAlignment = -4108
$ws.Range("N2:O2").VerticalAlignment = -4108

# N3:O3 "In "/"Out" labels use the bold look as well.
$ws.Range("N3:O3").Font.Name = "Aptos Narrow"
$ws.Range("N3:O3").Font.Size = 12
$ws.Range("N3:O3").Font.Bold = $true
$ws.Range("N3:O3").HorizontalAlignment = -4108
$ws.Range("N3:O3").VerticalAlignment = -4108

# Data cells K3:L8 and N4:O9 -> plain 12pt, centered (matches the B3:I8 look).
$dataRanges = @("K3:L8", "N4:O9")
foreach ($addr in $dataRanges) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Aptos Narrow"
    $r.Font.Size = 12
    $r.Font.Bold = $false
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}

# Thin left-hand border on the new K:L block to set it apart visually.
foreach ($addr in @("K2", "L2", "K3", "L3", "K4", "L4", "K5", "L5")) {
    $r = $ws.Range($addr)
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(7).Weight = 2
}

# Borders around the main table headers / N:O block (matches existing box style).
foreach ($addr in @("N2:O2", "N3:O3", "N4:O9")) {
    $r = $ws.Range($addr)
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(7).Weight = 2
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(8).Weight = 2
    $r.Borders.Item(9).LineStyle = 1
    $r.Borders.Item(9).Weight = 2
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(10).Weight = 2
}

# ---------------------------------------------------------------------
# 5) Merge the N2:O2 header cell (like the other category headers).
# ---------------------------------------------------------------------
$ws.Range("N2:O2").Merge()

# ---------------------------------------------------------------------
# 6) Threaded comment on N4, mirroring the D3 note.
# ---------------------------------------------------------------------
$excel.UserName = "Astha Bista"
$ws.Range("N4").AddCommentThreaded("Assumed 0. But maybe find source.")

# ---------------------------------------------------------------------
# 7) Make Sheet1 the active sheet / selection.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("P3").Select()
